# Auto-generated Excel COM-interop script to apply the Alpha_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3464.1428
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3464.1428
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3464.1428
$ws.Range("N40").Value = -3814.1428

$ws.Range("H53").Value = 1321.619
$ws.Range("I53").Value = 1815.6666
$ws.Range("J53").Value = 951.0833
$ws.Range("K53").Value = 1815.6666
$ws.Range("L53").Value = 951.0833
$ws.Range("M53").Value = -1178.6666
$ws.Range("N53").Value = -2225.0833

$ws.Range("H63").Value = 45000
$ws.Range("I63").Value = 45000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 45000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -44376
$ws.Range("N63").Value = $null

$ws.Range("H66").Value = 45000
$ws.Range("I66").Value = 45000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 135000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -131880
$ws.Range("N66").Value = $null

$ws.Range("H98").Value = 1891.9584
$ws.Range("I98").Value = 1870.4
$ws.Range("J98").Value = 1999.75
$ws.Range("K98").Value = 1870.4
$ws.Range("L98").Value = 1999.75
$ws.Range("M98").Value = -372.4000000000001

$ws.Range("H113").Value = 6645.9287
$ws.Range("I113").Value = 6480.75
$ws.Range("J113").Value = 6866.1665
$ws.Range("K113").Value = 6480.75
$ws.Range("L113").Value = 6866.1665
$ws.Range("M113").Value = -3226.75
$ws.Range("N113").Value = -13374.1665

$ws.Range("H122").Value = 1891.9584
$ws.Range("I122").Value = 1870.4
$ws.Range("J122").Value = 1999.75
$ws.Range("K122").Value = 5611.200000000001
$ws.Range("L122").Value = 5999.25
$ws.Range("M122").Value = -3161.200000000001

$ws.Range("H131").Value = 2074.2222
$ws.Range("I131").Value = 1161.3334
$ws.Range("J131").Value = 3900
$ws.Range("K131").Value = 3484.0002
$ws.Range("L131").Value = 11700
$ws.Range("M131").Value = 1555.9998

$ws.Range("H132").Value = 48464.09
$ws.Range("I132").Value = 48464.09
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 145392.27
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -142862.27
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3180.5806
$ws.Range("I32").Value = 3220.4666
$ws.Range("J32").Value = 1984
$ws.Range("K32").Value = 3220.4666
$ws.Range("L32").Value = 1984
$ws.Range("M32").Value = -2933.4666
$ws.Range("N32").Value = -2558

$ws.Range("H45").Value = 1569.2667
$ws.Range("I45").Value = 1496.4166
$ws.Range("J45").Value = 1860.6666
$ws.Range("K45").Value = 1496.4166
$ws.Range("L45").Value = 1860.6666
$ws.Range("M45").Value = -1119.4166
$ws.Range("N45").Value = -2614.6666

$ws.Range("H74").Value = 15432431
$ws.Range("I74").Value = 9259759
$ws.Range("J74").Value = 27777776
$ws.Range("K74").Value = 9259759
$ws.Range("L74").Value = 27777776
$ws.Range("M74").Value = -9258885

$ws.Range("H77").Value = 15432431
$ws.Range("I77").Value = 9259759
$ws.Range("J77").Value = 27777776
$ws.Range("K77").Value = 46298795
$ws.Range("L77").Value = 138888880
$ws.Range("M77").Value = -46294427

$ws.Range("H108").Value = 44249.5
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 44249.5
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 44249.5
$ws.Range("N108").Value = -51929.5

$ws.Range("H122").Value = 1891.375
$ws.Range("I122").Value = 1226.8
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 3680.4
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -1230.4
$ws.Range("N122").Value = -13897

$ws.Range("H132").Value = 12503872
$ws.Range("I132").Value = 2440.8057
$ws.Range("J132").Value = 125016750
$ws.Range("K132").Value = 7322.4171
$ws.Range("L132").Value = 375050250
$ws.Range("M132").Value = -4792.4171
$ws.Range("N132").Value = -375055310

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 3366.8
$ws.Range("I3").Value = 3613
$ws.Range("J3").Value = 2997.5
$ws.Range("K3").Value = 3613
$ws.Range("L3").Value = 2997.5
$ws.Range("M3").Value = -3500

$ws.Range("H31").Value = 2913.1333
$ws.Range("I31").Value = 2607.7
$ws.Range("J31").Value = 3524
$ws.Range("K31").Value = 2607.7
$ws.Range("L31").Value = 3524
$ws.Range("M31").Value = -2312.7
$ws.Range("N31").Value = -4114

$ws.Range("H34").Value = 2913.1333
$ws.Range("I34").Value = 2607.7
$ws.Range("J34").Value = 3524
$ws.Range("K34").Value = 2607.7
$ws.Range("L34").Value = 3524
$ws.Range("M34").Value = -2405.7
$ws.Range("N34").Value = -3928

$ws.Range("H92").Value = 87814
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 87814
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 87814
$ws.Range("N92").Value = -92806

$ws.Range("H96").Value = 47552
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 47552
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 47552
$ws.Range("N96").Value = -53044

$ws.Range("H114").Value = 44973
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 44973
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 44973
$ws.Range("N114").Value = -53651

$ws.Range("H122").Value = 44332.668
$ws.Range("I122").Value = 3199.2
$ws.Range("J122").Value = 250000
$ws.Range("K122").Value = 9597.599999999999
$ws.Range("L122").Value = 750000
$ws.Range("M122").Value = -7147.599999999999
$ws.Range("N122").Value = -754900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 547
$ws.Range("I12").Value = 249.25
$ws.Range("J12").Value = 695.875
$ws.Range("K12").Value = 747.75
$ws.Range("L12").Value = 2087.625
$ws.Range("M12").Value = -574.75
$ws.Range("N12").Value = -2433.625

$ws.Range("H17").Value = 433.33334
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 433.33334
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1300.00002
$ws.Range("M17").Value = $null
$ws.Range("N17").Value = -1638.00002

$ws.Range("H39").Value = 6650.875
$ws.Range("I39").Value = 3650.5
$ws.Range("J39").Value = 8451.1
$ws.Range("K39").Value = 10951.5
$ws.Range("L39").Value = 25353.3
$ws.Range("M39").Value = -10657.5
$ws.Range("N39").Value = -25941.3

$ws.Range("H69").Value = 2430
$ws.Range("I69").Value = 2400
$ws.Range("J69").Value = 2445
$ws.Range("K69").Value = 7200
$ws.Range("L69").Value = 7335
$ws.Range("M69").Value = -6389
$ws.Range("N69").Value = -8957

$ws.Range("H72").Value = 2430
$ws.Range("I72").Value = 2400
$ws.Range("J72").Value = 2445
$ws.Range("K72").Value = 21600
$ws.Range("L72").Value = 22005
$ws.Range("M72").Value = -17544
$ws.Range("N72").Value = -30117

$ws.Range("H80").Value = 3283.1667
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 3739.8
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 11219.4
$ws.Range("M80").Value = -2064
$ws.Range("N80").Value = -13091.4

$ws.Range("H83").Value = 3283.1667
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 3739.8
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 33658.2
$ws.Range("M83").Value = -4320
$ws.Range("N83").Value = -43018.2

$ws.Range("H114").Value = 2841
$ws.Range("I114").Value = 1446.5
$ws.Range("J114").Value = 3040.2144
$ws.Range("K114").Value = 4339.5
$ws.Range("L114").Value = 9120.643199999999
$ws.Range("M114").Value = -1085.5
$ws.Range("N114").Value = -15628.6432

$ws.Range("H131").Value = 350833.03
$ws.Range("I131").Value = 614.5714
$ws.Range("J131").Value = 405311.47
$ws.Range("K131").Value = 1843.7142
$ws.Range("L131").Value = 1215934.41
$ws.Range("M131").Value = 3196.2858

$ws.Range("H133").Value = 6775
$ws.Range("I133").Value = 6775
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 20325
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -15265

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 100000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 100000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 100000
$ws.Range("N63").Value = -101372

$ws.Range("H66").Value = 100000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 100000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 300000
$ws.Range("N66").Value = -306864

$ws.Range("H80").Value = 2481.1667
$ws.Range("I80").Value = 2869.6667
$ws.Range("J80").Value = 2092.6667
$ws.Range("K80").Value = 2869.6667
$ws.Range("L80").Value = 2092.6667
$ws.Range("M80").Value = -1871.6667
$ws.Range("N80").Value = -4088.6667

$ws.Range("H83").Value = 2481.1667
$ws.Range("I83").Value = 2869.6667
$ws.Range("J83").Value = 2092.6667
$ws.Range("K83").Value = 14348.3335
$ws.Range("L83").Value = 10463.3335
$ws.Range("M83").Value = -9356.333500000001
$ws.Range("N83").Value = -20447.3335

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2075.5
$ws.Range("I46").Value = 1195.5
$ws.Range("J46").Value = 2163.5
$ws.Range("K46").Value = 1195.5
$ws.Range("L46").Value = 2163.5
$ws.Range("M46").Value = -1007.5
$ws.Range("N46").Value = -2539.5

$ws.Range("H68").Value = 4555.7617
$ws.Range("I68").Value = 2741.2222
$ws.Range("J68").Value = 5916.6665
$ws.Range("K68").Value = 2741.2222
$ws.Range("L68").Value = 5916.6665
$ws.Range("M68").Value = -1992.2222
$ws.Range("N68").Value = -7414.6665

$ws.Range("H71").Value = 4555.7617
$ws.Range("I71").Value = 2741.2222
$ws.Range("J71").Value = 5916.6665
$ws.Range("K71").Value = 13706.111
$ws.Range("L71").Value = 29583.3325
$ws.Range("M71").Value = -9962.111000000001
$ws.Range("N71").Value = -37071.3325

$ws.Range("H122").Value = 3080.4285
$ws.Range("I122").Value = 3009.8462
$ws.Range("J122").Value = 3998
$ws.Range("K122").Value = 9029.5386
$ws.Range("L122").Value = 11994
$ws.Range("M122").Value = -6579.5386
$ws.Range("N122").Value = -16894

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 55028.5
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 55028.5
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 55028.5
$ws.Range("N64").Value = -55524.5

$ws.Range("H67").Value = 55028.5
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 55028.5
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 55028.5
$ws.Range("N67").Value = -56744.5

$ws.Range("H122").Value = 2557
$ws.Range("I122").Value = 2599.8333
$ws.Range("J122").Value = 2300
$ws.Range("K122").Value = 7799.499899999999
$ws.Range("L122").Value = 6900
$ws.Range("M122").Value = -5349.499899999999
